# Update the "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column after "Week" (shifts ASIN..is_holiday_week right by one)
#  - change Week labels from zero-padded (W01..W16) to unpadded (W1..W16)
#  - change is_holiday_week from numeric 0 to boolean FALSE

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1) Insert a new column B ("Week_Start_Date"); existing B..I (ASIN..is_holiday_week) shift to C..J
$ws.Columns.Item(2).Insert()

# 2) Header for the new column
$ws.Range("B1").Value = "Week_Start_Date"

# Week number (column A) -> Week start date (column B) mapping for the 16 forecast rows
$weekDates = @{
    2  = "2025-01-05"
    3  = "2025-01-12"
    4  = "2025-01-19"
    5  = "2025-01-26"
    6  = "2025-02-02"
    7  = "2025-02-09"
    8  = "2025-02-16"
    9  = "2025-02-23"
    10 = "2025-03-02"
    11 = "2025-03-09"
    12 = "2025-03-16"
    13 = "2025-03-23"
    14 = "2025-03-30"
    15 = "2025-04-06"
    16 = "2025-04-13"
    17 = "2025-04-20"
}

for ($row = 2; $row -le 17; $row++) {
    # Week label: drop the leading zero (W01 -> W1, ..., W16 stays W16)
    $weekNum = $row - 1
    $ws.Range("A$row").Value = "W$weekNum"

    # Week_Start_Date as literal text (not an Excel date serial number):
    # write it through a formula returning a text literal, then freeze the
    # result back into a plain value via copy / paste-values.
    $dateCell = $ws.Range("B$row")
    $dateCell.Formula = "=""" + $weekDates[$row] + """"
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)  # xlPasteValues

    # is_holiday_week (now column J) becomes a boolean instead of a number
    $ws.Range("J$row").Value = $false
}

$excel.CutCopyMode = 0
